$d = $word.ActiveDocument

# --- 1. Turn the "Table <SEQ>" caption field into a plain, static number ---
# The caption paragraph currently reads "Table 3: Sex-specific parameter
# values ..." where "3" is produced by a { SEQ Table \* ARABIC } field.
# The edit bumps the displayed number by one (3 -> 4) and bakes it in as
# literal text (no more field codes).
if ($d.Fields.Count -ge 1) {
    $field = $d.Fields.Item(1)
    $oldText = $field.Result.Text
    $oldNum = [int]$oldText
    $newText = [string]($oldNum + 1)

    # Unlink collapses the { begin / instrText / separate / result / end }
    # run sequence down to a single plain-text run holding the field's
    # current result, while leaving the surrounding runs untouched.
    $field.Unlink()

    # Locate that now-static result text within its paragraph and replace
    # it with the incremented value, keeping it as its own run.
    $para = $d.Paragraphs.Item(1).Range
    $hit = $para.Duplicate
    $found = $hit.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $hit.Text = $newText
        $afterNumber = $hit.End
    } else {
        $afterNumber = $para.Start
    }

    # --- 2. Relocate the "_GoBack" bookmark ---
    # It used to sit at the very end of the caption paragraph; it now sits
    # immediately after the table number (and before the trailing
    # ": Sex-specific ..." text).
    if ($d.Bookmarks.Exists("_GoBack")) {
        $goBack = $d.Bookmarks.Item("_GoBack")
        $goBack.Delete()
    }
    $bmRange = $d.Range($afterNumber, $afterNumber)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
